$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = '70.093.13'
$ws.Cells.Item(2, 5).Value = '  -1.05%  '
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = '3.591.89'
$ws.Cells.Item(3, 5).Value = '  -1.70%  '
$ws.Cells.Item(4, 5).Value = '  +0.13%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '575.71'
$ws.Cells.Item(5, 5).Value = '  -2.98%  '
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '188.15'
$ws.Cells.Item(6, 5).Value = '  -3.57%  '
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '3.587.99'
$ws.Cells.Item(7, 5).Value = '  -1.62%  '
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.628'
$ws.Cells.Item(8, 5).Value = '  -3.50%  '
$ws.Cells.Item(9, 5).Value = '  +0.02%  '
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '0.187'
$ws.Cells.Item(10, 5).Value = '  +2.87%  '
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.654'
$ws.Cells.Item(11, 5).Value = '  -3.45%  '
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '55.60'
$ws.Cells.Item(12, 5).Value = '  -5.22%  '
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '0.0000310'
$ws.Cells.Item(13, 5).Value = '  +5.36%  '
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '9.62'
$ws.Cells.Item(14, 5).Value = '  -3.66%  '
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '4.176.44'
$ws.Cells.Item(15, 5).Value = '  -1.38%  '
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '19.72'
$ws.Cells.Item(16, 5).Value = '  -1.24%  '
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '3.602.62'
$ws.Cells.Item(17, 5).Value = '  -1.42%  '
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '70.153.74'
$ws.Cells.Item(18, 5).Value = '  -0.91%  '
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '12.61'
$ws.Cells.Item(19, 5).Value = '  -1.71%  '
$ws.Cells.Item(20, 5).Value = '  -0.25%  '
$ws.Cells.Item(21, 5).Value = '  -3.27%  '
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '491.95'
$ws.Cells.Item(22, 5).Value = '  -0.05%  '
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '19.48'
$ws.Cells.Item(23, 5).Value = '  +2.92%  '
$ws.Cells.Item(24, 5).Value = '  -8.16%  '
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '95.77'
$ws.Cells.Item(25, 5).Value = '  +4.15%  '
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '4.34'
$ws.Cells.Item(26, 5).Value = '  -3.64%  '
$ws.Cells.Item(27, 2).Value = 'RenderToken'
$ws.Cells.Item(27, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '11.21'
$ws.Cells.Item(27, 5).Value = '  -2.50%  '
$ws.Cells.Item(28, 2).Value = 'ImmutableX'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '2.97'
$ws.Cells.Item(28, 5).Value = '  -6.64%  '
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '9.33'
$ws.Cells.Item(29, 5).Value = '  -3.11%  '
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '31.99'
$ws.Cells.Item(30, 5).Value = '  -3.46%  '
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '7.62'
$ws.Cells.Item(31, 5).Value = '  -3.62%  '
$ws.Cells.Item(32, 2).Value = 'OKB'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '67.12'
$ws.Cells.Item(32, 5).Value = '  +2.05%  '
$ws.Cells.Item(33, 2).Value = 'Cosmos'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '12.14'
$ws.Cells.Item(33, 5).Value = '  -1.66%  '
$ws.Cells.Item(34, 5).Value = '  -3.98%  '
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '574.06'
$ws.Cells.Item(35, 5).Value = '  -9.13%  '
$ws.Cells.Item(36, 2).Value = 'Fetch.AI'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '3.14'
$ws.Cells.Item(36, 5).Value = '  +9.47%  '
$ws.Cells.Item(37, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '38.42'
$ws.Cells.Item(37, 5).Value = '  -5.63%  '
$ws.Cells.Item(38, 2).Value = 'PEPE'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '0.0₃0807'
$ws.Cells.Item(38, 5).Value = '  -4.08%  '
$ws.Cells.Item(39, 2).Value = 'Dai'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '0.999'
$ws.Cells.Item(39, 5).Value = '  -0.09%  '
$ws.Cells.Item(40, 2).Value = 'dogwifhat'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '3.34'
$ws.Cells.Item(40, 5).Value = '  +13.29%  '
$ws.Cells.Item(41, 2).Value = 'TheGraph'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '0.395'
$ws.Cells.Item(41, 5).Value = '  -4.66%  '
$ws.Cells.Item(42, 5).Value = '  -1.14%  '
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '0.137'
$ws.Cells.Item(43, 5).Value = '  -7.44%  '
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '3.01'
$ws.Cells.Item(44, 5).Value = '  -5.37%  '
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '3.230.44'
$ws.Cells.Item(45, 5).Value = '  -2.73%  '
$ws.Cells.Item(46, 2).Value = 'VeChain'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '0.0438'
$ws.Cells.Item(46, 5).Value = '  -3.80%  '
$ws.Cells.Item(47, 2).Value = 'ApeXProtocol'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '3.43'
$ws.Cells.Item(47, 5).Value = '  +3.74%  '
$ws.Cells.Item(48, 2).Value = 'THORChain'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '9.69'
$ws.Cells.Item(48, 5).Value = '  +4.71%  '
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '0.137'
$ws.Cells.Item(49, 5).Value = '  -1.37%  '
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '1.00'
$ws.Cells.Item(50, 5).Value = '  +0.10%  '
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '3.17'
$ws.Cells.Item(51, 5).Value = '  -4.76%  '
